$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($Sheet, $Row, $Col, $Val)
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextCell $ws 2 4 "97.292.49"
Set-TextCell $ws 2 5 "  +5.18%  "

Set-TextCell $ws 3 4 "3.116.92"
Set-TextCell $ws 3 5 "  +0.30%  "

Set-TextCell $ws 4 5 "  -0.04%  "

Set-TextCell $ws 5 4 "238.32"
Set-TextCell $ws 5 5 "  +1.59%  "

Set-TextCell $ws 6 4 "610.31"
Set-TextCell $ws 6 5 "  -0.30%  "

Set-TextCell $ws 7 4 "1.11"
Set-TextCell $ws 7 5 "  +2.26%  "

Set-TextCell $ws 8 5 "  -0.77%  "

Set-TextCell $ws 9 5 "  +0.06%  "

Set-TextCell $ws 10 4 "3.113.62"
Set-TextCell $ws 10 5 "  +0.29%  "

Set-TextCell $ws 11 4 "0.786"
Set-TextCell $ws 11 5 "  -0.58%  "

Set-TextCell $ws 12 5 "  -0.30%  "

Set-TextCell $ws 13 4 "96.872.43"
Set-TextCell $ws 13 5 "  +4.97%  "

Set-TextCell $ws 14 5 "  -0.95%  "

Set-TextCell $ws 15 4 "33.88"
Set-TextCell $ws 15 5 "  +0.17%  "

Set-TextCell $ws 16 4 "5.42"
Set-TextCell $ws 16 5 "  +0.28%  "

Set-TextCell $ws 17 4 "3.693.48"
Set-TextCell $ws 17 5 "  +0.06%  "

Set-TextCell $ws 18 4 "3.113.38"
Set-TextCell $ws 18 5 "  +0.56%  "

Set-TextCell $ws 19 4 "3.55"
Set-TextCell $ws 19 5 "  -6.40%  "

Set-TextCell $ws 20 4 "526.66"
Set-TextCell $ws 20 5 "  +20.12%  "

Set-TextCell $ws 21 4 "14.53"
Set-TextCell $ws 21 5 "  +0.19%  "

Set-TextCell $ws 22 4 "5.66"
Set-TextCell $ws 22 5 "  -2.93%  "

Set-TextCell $ws 23 4 "0.0000194"
Set-TextCell $ws 23 5 "  -4.10%  "

Set-TextCell $ws 24 4 "8.82"
Set-TextCell $ws 24 5 "  -3.17%  "

Set-TextCell $ws 25 4 "5.49"
Set-TextCell $ws 25 5 "  -1.42%  "

Set-TextCell $ws 26 4 "87.78"
Set-TextCell $ws 26 5 "  +2.72%  "

Set-TextCell $ws 27 4 "11.56"
Set-TextCell $ws 27 5 "  +0.78%  "

Set-TextCell $ws 28 4 "3.281.10"
Set-TextCell $ws 28 5 "  +0.25%  "

Set-TextCell $ws 29 5 "  -0.10%  "

Set-TextCell $ws 30 4 "0.237"
Set-TextCell $ws 30 5 "  +0.97%  "

Set-TextCell $ws 31 5 "  -3.44%  "

Set-TextCell $ws 32 4 "0.124"
Set-TextCell $ws 32 5 "  +0.82%  "

Set-TextCell $ws 33 2 "InternetComputer(DFINITY)"
Set-TextCell $ws 33 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell $ws 33 4 "8.99"
Set-TextCell $ws 33 5 "  -1.73%  "

Set-TextCell $ws 34 2 "EthereumClassic"
Set-TextCell $ws 34 3 "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextCell $ws 34 4 "26.64"
Set-TextCell $ws 34 5 "  +3.92%  "

Set-TextCell $ws 35 2 "Binance-PegBSC-USD"
Set-TextCell $ws 35 3 "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextCell $ws 35 4 "0.868"
Set-TextCell $ws 35 5 "  -16.22%  "

Set-TextCell $ws 36 5 "  -7.37%  "

Set-TextCell $ws 37 4 "7.26"
Set-TextCell $ws 37 5 "  -10.09%  "

Set-TextCell $ws 38 2 "PancakeSwap"
Set-TextCell $ws 38 3 "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws 38 4 "1.87"
Set-TextCell $ws 38 5 "  -0.73%  "

Set-TextCell $ws 39 2 "Bittensor"
Set-TextCell $ws 39 3 "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell $ws 39 4 "480.86"
Set-TextCell $ws 39 5 "  +3.84%  "

Set-TextCell $ws 40 5 "  +1.47%  "

Set-TextCell $ws 41 5 "  +2.37%  "

Set-TextCell $ws 42 4 "1.23"
Set-TextCell $ws 42 5 "  -3.59%  "

Set-TextCell $ws 43 4 "3.58"
Set-TextCell $ws 43 5 "  -9.71%  "

Set-TextCell $ws 44 2 "USDe"
Set-TextCell $ws 44 3 "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextCell $ws 44 4 "1.00"
Set-TextCell $ws 44 5 "  -0.01%  "

Set-TextCell $ws 45 2 "dogwifhat"
Set-TextCell $ws 45 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell $ws 45 4 "3.17"
Set-TextCell $ws 45 5 "  -3.79%  "

Set-TextCell $ws 46 4 "162.09"
Set-TextCell $ws 46 5 "  +1.39%  "

Set-TextCell $ws 47 4 "1.91"
Set-TextCell $ws 47 5 "  +4.60%  "

Set-TextCell $ws 48 5 "  +1.48%  "

Set-TextCell $ws 49 4 "44.45"
Set-TextCell $ws 49 5 "  +1.55%  "

Set-TextCell $ws 50 4 "4.45"
Set-TextCell $ws 50 5 "  +2.70%  "

Set-TextCell $ws 51 2 "FirstDigitalUSD"
Set-TextCell $ws 51 3 "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell $ws 51 4 "0.998"
Set-TextCell $ws 51 5 "  +0.00%  "
